$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 7, 9, 13, 15 to match repulled data
$ws.Range("F7").Value = -9
$ws.Range("F9").Value = -8
$ws.Range("F13").Value = -6
$ws.Range("F15").Value = -7
